# Insert a new "PopulationCen" column between the existing "PopulationPPB"
# column (B) and "PersonCrimeAmt" column (old C, now D), filling in the
# new column's header and values for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns C:F to D:G by inserting a new blank column at C.
$ws.Columns.Item(3).Insert()

# New column header.
$ws.Range("C1").Value = "PopulationCen"

# New column values (Census population figures), rows 2-10.
$values = @(585436, 595410, 604285, 609970, 620647, 631539, 643136, 648630, 643115)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

# Only the brand-new column needs an explicit width; the inserted column
# naturally pushed the old formatting of C:F over to D:G, preserving their
# existing (bestFit) widths.
$ws.Columns.Item(3).ColumnWidth = 12.166666666666666

# Update the selection to reflect the post-edit state.
$ws.Range("E21").Select()
